$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 26 - Wages Expense
$ws.Range("B26").Value = "Wages Expense"
$ws.Range("D26").Formula = "=60000+240000"

# Row 27 - TRANSFER BCA
$ws.Range("B27").Value = "TRANSFER BCA"
$ws.Range("D27").Formula = "=17610000+23880000+2295000+6410000+221000"

# Row 28 - TAX - IURAN ARIESTA (plain value, not a formula)
$ws.Range("B28").Value = "TAX - IURAN ARIESTA"
$ws.Range("D28").Value = 660000

# Row 29 - TAX - P.Tata
$ws.Range("B29").Value = "TAX - P.Tata"
$ws.Range("D29").Formula = "=200000"

# Row 30 - A/R
$ws.Range("B30").Value = "A/R"
$ws.Range("C30").Formula = "=6410000"

# Row 31 - SALES - cash/retail
$ws.Range("B31").Value = "SALES - cash/retail"
$ws.Range("C31").Formula = "=44283375+26614925"

# Row 32 - A/P
$ws.Range("B32").Value = "A/P"
$ws.Range("D32").Formula = "=2300000"

# Row 33 - SELISIH - lebih (plain value)
$ws.Range("B33").Value = "SELISIH - lebih"
$ws.Range("C33").Value = 89700

# Row 34 - SETOR KE BANK
$ws.Range("B34").Value = "SETOR KE BANK"
$ws.Range("D34").Formula = "=24000000"

# Row 35 - new date
$ws.Range("A35").Value = 44351

# Update selection to match end-of-day edit position
$null = $ws.Range("B35").Select()

$wb.Save()
